# Update column F (row 2 through 25) of Sheet1 with new computed values
# for the "380 kV" case, as described by the commit "case with 380 kV done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = [ordered]@{
    "F2"  = 0.4443680307746121
    "F3"  = 0.3878228170618172
    "F4"  = 0.3531389305169483
    "F5"  = 0.3390132514313251
    "F6"  = 0.336668177824194
    "F7"  = 0.3529483938345521
    "F8"  = 0.4248636149814047
    "F9"  = 0.5661985755041457
    "F10" = 0.6702781546542269
    "F11" = 0.7176906081379002
    "F12" = 0.7356546913071611
    "F13" = 0.7317853510981394
    "F14" = 0.7191683204515442
    "F15" = 0.7114413442032514
    "F16" = 0.6671810134426437
    "F17" = 0.6400460337125793
    "F18" = 0.6244449056556647
    "F19" = 0.619163680173358
    "F20" = 0.642933953830422
    "F21" = 0.7228739723492197
    "F22" = 0.7751780083420101
    "F23" = 0.7472568307830727
    "F24" = 0.6416283278901602
    "F25" = 0.5279251897347166
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
